$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for every data row (2..540)
# from 2023-09-21 (45190) to 2023-09-23 (45192).
$ws.Range("C2:C540").Value = 45192

# Row 4 additionally lost the "Motaggsvamp" signal species, which reduces
# the NT (J), Rödlistade (O) and Alla arter (Q) counts by one each, and
# removes the corresponding line from the species list (R).
$speciesList = $ws.Range("R4").Value()
$speciesList = $speciesList -replace "Motaggsvamp`r`n", ""
$ws.Range("R4").Value = $speciesList

$ws.Range("J4").Value = 13
$ws.Range("O4").Value = 17
$ws.Range("Q4").Value = 20

# Updating the wrapped R4 text can trigger an autofit row-height recalculation;
# restore the original fixed row height so only the intended cell values change.
$ws.Rows.Item(4).RowHeight = 15
